$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Pune"
$ws.Range("C1").Value = "Location"
